# Fruta / hortaliza, semanal
# Re-shuffle the Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) values across rows 2-20
# (skipping rows 3, 12 and 18 which are unaffected) according to the row
# mapping derived from the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that get rearranged, keyed by
# row number, before any writes happen (so that source data used later in
# the loop is never the already-modified value).
$cols = @("D", "M", "N", "O", "P", "S")
$rows = @(2, 4, 5, 6, 7, 8, 9, 10, 11, 13, 14, 15, 16, 17, 19, 20)

$original = @{}
foreach ($r in $rows) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Range("$c$r").Value2
    }
    $original[$r] = $vals
}

# destination row -> source row (values copied from the source row's
# original contents into the destination row)
$mapping = @{
    2  = 14
    4  = 13
    5  = 7
    6  = 2
    7  = 6
    8  = 17
    9  = 8
    10 = 19
    11 = 5
    13 = 15
    14 = 11
    15 = 16
    16 = 4
    17 = 20
    19 = 9
    20 = 10
}

foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    $srcVals = $original[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcVals[$c]
    }
}
